$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: repurpose as a compact "summary" row reusing the Maldives title
# and a new, shorter Portugal title (mirrors B1 without "на майские").
$ws.Range("A6").Value = "Серфинг на Мальдивах"
$ws.Range("B6").Value = "Серф-кемп в Португалии"

# Rows 7-13 held the long descriptive paragraphs; they are emptied out,
# leaving just their (now much shorter) row heights.
$ws.Range("A7:B13").ClearContents()

# Shrink rows 6-13 down to their new compact heights.
$ws.Rows.Item(6).RowHeight = 22.5
$ws.Rows.Item(7).RowHeight = 3
$ws.Rows.Item(8).RowHeight = 9
$ws.Rows.Item(9).RowHeight = 4.5
$ws.Rows.Item(10).RowHeight = 8.25
$ws.Rows.Item(11).RowHeight = 4.5
$ws.Rows.Item(12).RowHeight = 6.75
$ws.Rows.Item(13).RowHeight = 6.75

# Scroll view back to the top and move the live selection to B6.
$ws.Range("B6").Select()
